$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3 (Leve Item ID 18511)
$ws.Range("H3").Value = 27378.143
$ws.Range("J3").Value = 27378.143
$ws.Range("L3").Value = 27378.143
$ws.Range("N3").Value = -27606.143
# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 134.83333
$ws.Range("J9").Value = 53
$ws.Range("L9").Value = 53
$ws.Range("N9").Value = -391
# Row 29 (Leve Item ID 4575)
$ws.Range("H29").Value = 4196.077
$ws.Range("J29").Value = 4900
$ws.Range("L29").Value = 14700
$ws.Range("N29").Value = -15262
# Row 31 (Leve Item ID 4576)
$ws.Range("H31").Value = 87
$ws.Range("I31").Value = 87
$ws.Range("K31").Value = 261
$ws.Range("M31").Value = -31
# Row 38 (Leve Item ID 4599)
$ws.Range("H38").Value = 262
$ws.Range("I38").Value = 262
$ws.Range("K38").Value = 786
$ws.Range("M38").Value = -414
# Row 101 (Leve Item ID 19884)
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 102 (Leve Item ID 18511)
$ws.Range("H102").Value = 27378.143
$ws.Range("J102").Value = 27378.143
$ws.Range("L102").Value = 27378.143
$ws.Range("N102").Value = -33868.143
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 9997.044
$ws.Range("I132").Value = 9997.044
$ws.Range("K132").Value = 29991.132
$ws.Range("M132").Value = -27461.132

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4 (Leve Item ID 5071)
$ws.Range("H4").Value = 227
$ws.Range("J4").Value = 96.666664
$ws.Range("L4").Value = 96.666664
$ws.Range("N4").Value = -328.666664
# Row 5 (Leve Item ID 5091)
$ws.Range("H5").Value = 149
$ws.Range("I5").Value = 149
$ws.Range("K5").Value = 149
$ws.Range("M5").Value = -37
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 4762820.5
$ws.Range("I32").Value = 961.5
$ws.Range("K32").Value = 961.5
$ws.Range("M32").Value = -674.5
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 2815.923
$ws.Range("I45").Value = 2176.3
$ws.Range("K45").Value = 2176.3
$ws.Range("M45").Value = -1799.3
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 2150.7058
$ws.Range("I74").Value = 1823.7333
$ws.Range("J74").Value = 4603
$ws.Range("K74").Value = 1823.7333
$ws.Range("L74").Value = 4603
$ws.Range("M74").Value = -949.7333000000001
$ws.Range("N74").Value = -6351
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 2150.7058
$ws.Range("I77").Value = 1823.7333
$ws.Range("J77").Value = 4603
$ws.Range("K77").Value = 9118.666500000001
$ws.Range("L77").Value = 23015
$ws.Range("M77").Value = -4750.666500000001
$ws.Range("N77").Value = -31751

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4 (Leve Item ID 5091)
$ws.Range("H4").Value = 149
$ws.Range("I4").Value = 149
$ws.Range("K4").Value = 149
$ws.Range("M4").Value = -34
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 35718960
$ws.Range("I107").Value = 55557270
$ws.Range("K107").Value = 55557270
$ws.Range("M107").Value = -55555350

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 25 (Leve Item ID 1895)
$ws.Range("H25").Value = 2705.5
$ws.Range("I25").Value = 2705.5
$ws.Range("K25").Value = 2705.5
$ws.Range("M25").Value = -2531.5
# Row 28 (Leve Item ID 18348)
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 5835.8965
$ws.Range("I31").Value = 2925.9412
$ws.Range("K31").Value = 2925.9412
$ws.Range("M31").Value = -2630.9412
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 5835.8965
$ws.Range("I34").Value = 2925.9412
$ws.Range("K34").Value = 2925.9412
$ws.Range("M34").Value = -2723.9412
# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 1220.2667
$ws.Range("I132").Value = 1220.2667
$ws.Range("K132").Value = 3660.800099999999
$ws.Range("M132").Value = -1130.800099999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 1455.3636
$ws.Range("I113").Value = 966
$ws.Range("J113").Value = 1638.875
$ws.Range("K113").Value = 2898
$ws.Range("L113").Value = 4916.625
$ws.Range("M113").Value = -728
$ws.Range("N113").Value = -9256.625
# Row 116 (Leve Item ID 27866)
$ws.Range("H116").Value = 1485.4
$ws.Range("I116").Value = 1643
$ws.Range("J116").Value = 1249
$ws.Range("K116").Value = 4929
$ws.Range("L116").Value = 3747
$ws.Range("M116").Value = -1487
$ws.Range("N116").Value = -10631

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 1393.5
$ws.Range("I70").Value = 1393.5
$ws.Range("K70").Value = 1393.5
$ws.Range("M70").Value = -1123.5
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 1393.5
$ws.Range("I73").Value = 1393.5
$ws.Range("K73").Value = 1393.5
$ws.Range("M73").Value = -457.5
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 2856.1667
$ws.Range("I80").Value = 2661.2727
$ws.Range("K80").Value = 2661.2727
$ws.Range("M80").Value = -1663.2727
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 2856.1667
$ws.Range("I83").Value = 2661.2727
$ws.Range("K83").Value = 13306.3635
$ws.Range("M83").Value = -8314.363499999999
# Row 101 (Leve Item ID 18513)
$ws.Range("H101").Value = 45331.332
$ws.Range("J101").Value = 45331.332
$ws.Range("L101").Value = 45331.332
$ws.Range("N101").Value = -51821.332
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 57967.832
$ws.Range("I132").Value = 64838.812
$ws.Range("K132").Value = 194516.436
$ws.Range("M132").Value = -191986.436

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 9 (Leve Item ID 1685)
$ws.Range("H9").Value = 137.5
$ws.Range("I9").Value = 137.5
$ws.Range("K9").Value = 137.5
$ws.Range("M9").Value = 86.5
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 1367.3529
$ws.Range("I55").Value = 1088.6666
$ws.Range("J55").Value = 1680.875
$ws.Range("K55").Value = 1088.6666
$ws.Range("L55").Value = 1680.875
$ws.Range("M55").Value = -915.6666
$ws.Range("N55").Value = -2026.875
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3899
$ws.Range("I132").Value = 3899
$ws.Range("K132").Value = 11697
$ws.Range("M132").Value = -9167

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 3 (Leve Item ID 3309)
$ws.Range("H3").Value = 3499.6667
$ws.Range("I3").Value = 4499.5
$ws.Range("K3").Value = 4499.5
$ws.Range("M3").Value = -4385.5
# Row 4 (Leve Item ID 2996)
$ws.Range("H4").Value = 17028.572
$ws.Range("I4").Value = 52750
$ws.Range("J4").Value = 2740
$ws.Range("K4").Value = 52750
$ws.Range("L4").Value = 2740
$ws.Range("M4").Value = -52637
$ws.Range("N4").Value = -2966
# Row 5 (Leve Item ID 3515)
$ws.Range("H5").Value = 7676250
$ws.Range("J5").Value = 3859285.8
$ws.Range("L5").Value = 3859285.8
$ws.Range("N5").Value = -3859509.8
# Row 6 (Leve Item ID 3000)
$ws.Range("H6").Value = 499.5
$ws.Range("I6").Value = 499.5
$ws.Range("K6").Value = 499.5
$ws.Range("M6").Value = -384.5
# Row 70 (Leve Item ID 11979)
$ws.Range("H70").Value = 36666.332
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630
# Row 73 (Leve Item ID 11979)
$ws.Range("H73").Value = 36666.332
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184
# Row 80 (Leve Item ID 10911)
$ws.Range("H80").Value = 71150
$ws.Range("J80").Value = 71150
$ws.Range("L80").Value = 71150
$ws.Range("N80").Value = -73146
# Row 83 (Leve Item ID 10911)
$ws.Range("H83").Value = 71150
$ws.Range("J83").Value = 71150
$ws.Range("L83").Value = 213450
$ws.Range("N83").Value = -223434
# Row 95 (Leve Item ID 18243)
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 781.9375
$ws.Range("I132").Value = 781.9375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2345.8125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 184.1875
$ws.Range("N132").ClearContents()
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 3141.5
$ws.Range("I136").Value = 1985.5
$ws.Range("J136").Value = 4682.8335
$ws.Range("K136").Value = 5956.5
$ws.Range("L136").Value = 14048.5005
$ws.Range("M136").Value = -3406.5
$ws.Range("N136").Value = -19148.5005
